# Apply edits for BarycentricCoordinates workbook update:
# - Reduce grid square size (H2, H3) from 10 to 5
# - Update a handful of existing lookup-table cells (J9:L11) with refined data
# - Add a new, more detailed height sample grid in columns U:AD (rows 2-11)
# - Update the active selection to AC22

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Core parameters ---
$ws.Range("H2").Value = 5
$ws.Range("H3").Value = 5

# --- Refined values in the existing J:S "Heights" lookup table ---
$ws.Range("J9").Value = 1.9149378500000001
$ws.Range("K9").Value = 2.1694214299999999
$ws.Range("L9").Value = 3
$ws.Range("J10").Value = 1.1240209299999999
$ws.Range("K10").Value = 1.6814620499999999
$ws.Range("L10").Value = 1.70000017
$ws.Range("K11").Value = 0.5
$ws.Range("L11").Value = 1

# --- New supplementary height-sample grid, columns U:AD, rows 2-11 ---
$ws.Range("U2").Value = 7.7241400000000002
$ws.Range("V2").Value = 8.2077899999999993
$ws.Range("W2").Value = 9.1060199999999991
$ws.Range("X2").Value = 9.24
$ws.Range("Y2").Value = 9.7172999999999998
$ws.Range("Z2").Value = 10.6394
$ws.Range("AA2").Value = 11
$ws.Range("AB2").Value = 10.535500000000001
$ws.Range("AC2").Value = 10.081200000000001
$ws.Range("AD2").Value = 9.8813600000000008
$ws.Range("U3").Value = 7.1181799999999997
$ws.Range("V3").Value = 8.3627500000000001
$ws.Range("W3").Value = 8.8815799999999996
$ws.Range("X3").Value = 9.0268300000000004
$ws.Range("Y3").Value = 9.2234300000000005
$ws.Range("Z3").Value = 10.033099999999999
$ws.Range("AA3").Value = 10.4359
$ws.Range("AB3").Value = 10.0405
$ws.Range("AC3").Value = 9.9442500000000003
$ws.Range("AD3").Value = 9.7623899999999999
$ws.Range("U4").Value = 6.4328900000000004
$ws.Range("V4").Value = 7.6282100000000002
$ws.Range("W4").Value = 8.8015899999999991
$ws.Range("X4").Value = 8.4420300000000008
$ws.Range("Y4").Value = 7.9097400000000002
$ws.Range("Z4").Value = 8.3754399999999993
$ws.Range("AA4").Value = 8.5328999999999997
$ws.Range("AB4").Value = 9.6984100000000009
$ws.Range("AC4").Value = 9.4666700000000006
$ws.Range("AD4").Value = 8.6892700000000005
$ws.Range("U5").Value = 5.15306
$ws.Range("V5").Value = 6.0157499999999997
$ws.Range("W5").Value = 7.4549200000000004
$ws.Range("X5").Value = 7.2142900000000001
$ws.Range("Y5").Value = 7.8365499999999999
$ws.Range("Z5").Value = 7.8296000000000001
$ws.Range("AA5").Value = 7.1428599999999998
$ws.Range("AB5").Value = 7.3333300000000001
$ws.Range("AC5").Value = 7.6296299999999997
$ws.Range("AD5").Value = 7.7142900000000001
$ws.Range("U6").Value = 4.1310799999999999
$ws.Range("V6").Value = 4.7524300000000004
$ws.Range("W6").Value = 5.3337199999999996
$ws.Range("X6").Value = 5.3390399999999998
$ws.Range("Y6").Value = 6.5378800000000004
$ws.Range("Z6").Value = 7.2430599999999998
$ws.Range("AA6").Value = 6.96
$ws.Range("AB6").Value = 6.7254899999999997
$ws.Range("AC6").Value = 6.4444499999999998
$ws.Range("AD6").Value = 6.93398
$ws.Range("U7").Value = 3.9420600000000001
$ws.Range("V7").Value = 3.8097400000000001
$ws.Range("W7").Value = 4.2846299999999999
$ws.Range("X7").Value = 5.0022000000000002
$ws.Range("Y7").Value = 5.2476200000000004
$ws.Range("Z7").Value = 4.6969700000000003
$ws.Range("AA7").Value = 4.6478000000000002
$ws.Range("AB7").Value = 5.7160500000000001
$ws.Range("AC7").Value = 5.88889
$ws.Range("AD7").Value = 5.3121
$ws.Range("U8").Value = 3.23529
$ws.Range("V8").Value = 2.8014199999999998
$ws.Range("W8").Value = 2.9716300000000002
$ws.Range("X8").Value = 3.5714299999999999
$ws.Range("Y8").Value = 3.5628099999999998
$ws.Range("Z8").Value = 3.6330900000000002
$ws.Range("AA8").Value = 3
$ws.Range("AB8").Value = 4.2415900000000004
$ws.Range("AC8").Value = 3.9317299999999999
$ws.Range("AD8").Value = 3.2857099999999999
$ws.Range("U9").Value = 2.1865800000000002
$ws.Range("V9").Value = 1.3596699999999999
$ws.Range("W9").Value = 1.3777299999999999
$ws.Range("X9").Value = 1.6666700000000001
$ws.Range("Y9").Value = 1.2795700000000001
$ws.Range("Z9").Value = 2.3867400000000001
$ws.Range("AA9").Value = 2.90055
$ws.Range("AB9").Value = 3.2544499999999998
$ws.Range("AC9").Value = 2.7066699999999999
$ws.Range("AD9").Value = 2.1352899999999999
$ws.Range("U10").Value = 1.8192600000000001
$ws.Range("V10").Value = 1.0701799999999999
$ws.Range("W10").Value = 1.0740700000000001
$ws.Range("X10").Value = 1.0238700000000001
$ws.Range("Y10").Value = 0.49019000000000001
$ws.Range("Z10").Value = 0.62963000000000002
$ws.Range("AA10").Value = 1.4512100000000001
$ws.Range("AB10").Value = 1.8676600000000001
$ws.Range("AC10").Value = 2.44279
$ws.Range("AD10").Value = 2.0850900000000001
$ws.Range("U11").Value = 1.30769
$ws.Range("V11").Value = 1.1126199999999999
$ws.Range("W11").Value = 1.0465100000000001
$ws.Range("X11").Value = 0.53846000000000005
$ws.Range("Y11").Value = 0.08306
$ws.Range("Z11").Value = 0.0885
$ws.Range("AA11").Value = 0.79310000000000003
$ws.Range("AB11").Value = 1.4761899999999999
$ws.Range("AC11").Value = 1.81481
$ws.Range("AD11").Value = 1.9213499999999999

# --- Restore the active cell selection ---
$ws.Range("AC22").Select()
